# Apply the data corrections to the "annotation" sheet (bounding-box /
# category values for a set of traffic-sign images) and move the
# selection to reflect where the author was last working in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- row 2 (img (1).png) ---
$ws.Range("B2").Value = 109
$ws.Range("C2").Value = 90
$ws.Range("D2").Value = 19
$ws.Range("E2").Value = 14
$ws.Range("F2").Value = 85
$ws.Range("G2").Value = 79

# --- row 11 ---
$ws.Range("B11").Value = 124
$ws.Range("C11").Value = 113
$ws.Range("D11").Value = 27
$ws.Range("E11").Value = 28
$ws.Range("F11").Value = 89
$ws.Range("G11").Value = 90
$ws.Range("H11").Value = 5

# --- row 20 ---
$ws.Range("B20").Value = 109
$ws.Range("C20").Value = 98
$ws.Range("D20").Value = 19
$ws.Range("E20").Value = 11
$ws.Range("F20").Value = 91
$ws.Range("G20").Value = 83
$ws.Range("H20").Value = 5

# --- row 25 ---
$ws.Range("B25").Value = 59
$ws.Range("C25").Value = 61
$ws.Range("D25").Value = 5
$ws.Range("E25").Value = 8
$ws.Range("F25").Value = 53
$ws.Range("G25").Value = 55
$ws.Range("H25").Value = 4

# --- row 27 ---
$ws.Range("B27").Value = 110
$ws.Range("C27").Value = 99
$ws.Range("D27").Value = 26
$ws.Range("E27").Value = 18
$ws.Range("F27").Value = 79
$ws.Range("G27").Value = 74
$ws.Range("H27").Value = 5

# --- row 29 ---
$ws.Range("B29").Value = 118
$ws.Range("C29").Value = 109
$ws.Range("D29").Value = 12
$ws.Range("E29").Value = 6
$ws.Range("F29").Value = 109
$ws.Range("G29").Value = 97
$ws.Range("H29").Value = 5

# --- row 30 (B30 unchanged) ---
$ws.Range("C30").Value = 123
$ws.Range("D30").Value = 19
$ws.Range("E30").Value = 20
$ws.Range("F30").Value = 90
$ws.Range("G30").Value = 90
$ws.Range("H30").Value = 6

# --- row 34 ---
$ws.Range("B34").Value = 114
$ws.Range("C34").Value = 114
$ws.Range("D34").Value = 22
$ws.Range("E34").Value = 21
$ws.Range("F34").Value = 97
$ws.Range("G34").Value = 95
$ws.Range("H34").Value = 7

# --- row 36 ---
$ws.Range("B36").Value = 90
$ws.Range("C36").Value = 84
$ws.Range("D36").Value = 21
$ws.Range("E36").Value = 18
$ws.Range("F36").Value = 64
$ws.Range("G36").Value = 61
$ws.Range("H36").Value = 22

# --- row 37 ---
$ws.Range("B37").Value = 112
$ws.Range("C37").Value = 107
$ws.Range("D37").Value = 16
$ws.Range("E37").Value = 14
$ws.Range("F37").Value = 96
$ws.Range("G37").Value = 89
$ws.Range("H37").Value = 26

# --- row 38 ---
$ws.Range("B38").Value = 96
$ws.Range("C38").Value = 94
$ws.Range("D38").Value = 19
$ws.Range("E38").Value = 1
$ws.Range("F38").Value = 73
$ws.Range("G38").Value = 70
$ws.Range("H38").Value = 35

# --- row 39 ---
$ws.Range("B39").Value = 83
$ws.Range("C39").Value = 86
$ws.Range("D39").Value = 15
$ws.Range("E39").Value = 18
$ws.Range("F39").Value = 63
$ws.Range("G39").Value = 60
$ws.Range("H39").Value = 39

# --- row 51 ---
$ws.Range("B51").Value = 145
$ws.Range("C51").Value = 134
$ws.Range("D51").Value = 32
$ws.Range("E51").Value = 26
$ws.Range("F51").Value = 115
$ws.Range("G51").Value = 106
$ws.Range("H51").Value = 42

# --- row 52 ---
$ws.Range("B52").Value = 128
$ws.Range("C52").Value = 118
$ws.Range("D52").Value = 33
$ws.Range("E52").Value = 27
$ws.Range("F52").Value = 101
$ws.Range("G52").Value = 91
$ws.Range("H52").Value = 43

# --- row 60 ---
$ws.Range("B60").Value = 110
$ws.Range("C60").Value = 99
$ws.Range("D60").Value = 17
$ws.Range("E60").Value = 19
$ws.Range("F60").Value = 85
$ws.Range("G60").Value = 80
$ws.Range("H60").Value = 47

# --- row 63 ---
$ws.Range("B63").Value = 92
$ws.Range("C63").Value = 87
$ws.Range("D63").Value = 16
$ws.Range("E63").Value = 15
$ws.Range("F63").Value = 72
$ws.Range("G63").Value = 68
$ws.Range("H63").Value = 3

# --- row 70 ---
$ws.Range("B70").Value = 195
$ws.Range("C70").Value = 185
$ws.Range("D70").Value = 38
$ws.Range("E70").Value = 24
$ws.Range("F70").Value = 164
$ws.Range("G70").Value = 154
$ws.Range("H70").Value = 26

# --- row 73 ---
$ws.Range("B73").Value = 200
$ws.Range("C73").Value = 187
$ws.Range("D73").Value = 32
$ws.Range("E73").Value = 13
$ws.Range("F73").Value = 169
$ws.Range("G73").Value = 176
$ws.Range("H73").Value = 6

# --- row 76 ---
$ws.Range("B76").Value = 103
$ws.Range("C76").Value = 97
$ws.Range("D76").Value = 12
$ws.Range("E76").Value = 13
$ws.Range("F76").Value = 92
$ws.Range("G76").Value = 88
$ws.Range("H76").Value = 7

# --- row 78 ---
$ws.Range("B78").Value = 113
$ws.Range("C78").Value = 106
$ws.Range("D78").Value = 20
$ws.Range("E78").Value = 17
$ws.Range("F78").Value = 98
$ws.Range("G78").Value = 90
$ws.Range("H78").Value = 7

# Move the view/selection to where the author ended up working (row 64
# scrolled into view, cell I69 selected) to match the saved workbook state.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 64
$ws.Range("I69").Select()
